$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riesgos")

# --- Minuta de la junta realizada el 13 de marzo de 2015 ---
# Nueva fila de riesgo (fila 12), previamente vacia.

$ws.Range("C12").Value = 42076
$ws.Range("D12").Value = "Retraso con las actividades establecidas por la asistencia al curso de blender."
$ws.Range("E12").Value = "Retraso del proyecto."
$ws.Range("F12").Value = "Alta"
$ws.Range("G12").Value = "Baja"
$ws.Range("H12").Formula = "=IF(F12=""Alta"",3,IF(F12=""Media"", 2, IF(F12=""Baja"",1, 0)))*IF(G12=""Alta"",3,IF(G12=""Media"", 2, IF(G12=""Baja"",1, 0)))"
$ws.Range("I12").Value = "Dedicarle horas extras al proyecto"
$ws.Range("J12").Value = "Cumplir con los tiempos establecidos"
$ws.Range("K12").Value = "JASS"

# La fila crece para acomodar el texto de la descripcion del riesgo.
$ws.Range("D12").EntireRow.RowHeight = 25.5

# Estado de la ventana/hoja al guardar: posicion y zoom usados durante la revision.
$ws.Range("C12").Select()
$excel.ActiveWindow.Zoom = 140
